# Adding test cases: Validate search by Filter on Search Skills page
#
# The ShareSkill sheet's first data row (row 2) is re-dated: the "Start
# date" (H2) moves from 6/10/2022 to 7/10/2022 and the "End date" (I2)
# moves from 7/31/2022 to 8/31/2022 (both stored as Excel date serials).
# The sheet's active/selected cell also moves to I3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start date / End date for the first ShareSkill row.
$ws.Range("H2").Value = 44752
$ws.Range("I2").Value = 44804

# Leave the cursor on the End date of the second row, matching the
# selection recorded for this sheet after the edit.
$ws.Range("I3").Select() | Out-Null
